$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 6 ("Game" row): Status -> Completed, Date of completion -> 6 November 2021
$statusCell = $t.Cell(6, 5)
$statusCell.Range.Text = "Completed"

$dateCell = $t.Cell(6, 6)
$dateCell.Range.Text = "6 November 2021"

# Clean up the "Niya & Tereza" cells that currently have a spell-check
# proofErr wrapper around "Niya" (rows 2, 4, 5 and 6 in the "Assign" column)
# by deleting the cell content and re-inserting plain text, which collapses
# the runs into a single run with no proofErr markers.
$assignRows = 2, 4, 5, 6
foreach ($rowIdx in $assignRows) {
    $cell = $t.Cell($rowIdx, 7)
    $range = $cell.Range
    $range.Delete()
    $range.InsertAfter("Niya & Tereza")
}
